# Trade #47 closed at 2026-02-17 21:07:39 - unknown UNKNOWN +0.000%
#
# - Closes trade #75 (early_exit) on both the "All Trades" and
#   "MarketMaking" sheets.
# - Appends a brand-new open trade #108 to both sheets.
# - Refreshes the roll-up numbers on "Summary" and "Strategy Status".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.51   # Current Capital
$summary.Range("B4").Value = 0.31      # Total P&L $
$summary.Range("B6").Value = 75        # Total Trades
$summary.Range("B7").Value = 34        # Winning Trades
$summary.Range("B9").Value = 45.33     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.51     # Capital
$status.Range("D5").Value = 42         # Trades
$status.Range("E5").Value = 0.2        # P&L $
$status.Range("F5").Value = 0.51       # P&L %
$status.Range("G5").Value = 47.62      # Win Rate %

# ---------------------------------------------------------------
# All Trades sheet
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry,
#          G Exit, H Status, I P&L%, J P&L$, K CapAfter, L ExitReason,
#          M Duration, N EntrySlip, O ExitSlip, P Confidence, Q EntryReason
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close trade #75 (row 76) - early exit
$allTrades.Range("G76").Value = 0.88
$allTrades.Range("H76").Value = "CLOSED"
$allTrades.Range("I76").Value = 2.3256
$allTrades.Range("J76").Value = 0.02
$allTrades.Range("K76").Value = 100.51
$allTrades.Range("L76").Value = "early_exit"
$allTrades.Range("M76").Value = 0.14

# Append new open trade #108 (row 109)
$allTrades.Range("A109").Value = 108
$allTrades.Range("B109").NumberFormat = "@"
$allTrades.Range("B109").Value = "2026-02-17"
$allTrades.Range("B109").Style = "Normal"
$allTrades.Range("C109").NumberFormat = "@"
$allTrades.Range("C109").Value = "21:07:32"
$allTrades.Range("C109").Style = "Normal"
$allTrades.Range("D109").Value = "MarketMaking"
$allTrades.Range("E109").Value = "UP"
$allTrades.Range("F109").Value = 0.86
$allTrades.Range("H109").Value = "OPEN"
$allTrades.Range("I109").Value = 0
$allTrades.Range("J109").Value = 0
$allTrades.Range("K109").Value = 100.4914872031006
$allTrades.Range("M109").Value = 0
$allTrades.Range("N109").Value = 0
$allTrades.Range("O109").Value = 0
$allTrades.Range("P109").Value = 0.6
$allTrades.Range("Q109").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------
# MarketMaking sheet
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry,
#          G Exit, H Status, I P&L%, J P&L$, K CapAfter, L EntrySlip,
#          M ExitSlip, N Confidence, O EntryReason, P ExitReason, Q Duration
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# Close trade #75 (row 43) - early exit
$marketMaking.Range("G43").Value = 0.88
$marketMaking.Range("H43").Value = "CLOSED"
$marketMaking.Range("I43").Value = 2.3256
$marketMaking.Range("J43").Value = 0.02
$marketMaking.Range("K43").Value = 100.51
$marketMaking.Range("P43").Value = "early_exit"
$marketMaking.Range("Q43").Value = 0.14

# Append new open trade #108 (row 76)
$marketMaking.Range("A76").Value = 108
$marketMaking.Range("B76").NumberFormat = "@"
$marketMaking.Range("B76").Value = "2026-02-17"
$marketMaking.Range("B76").Style = "Normal"
$marketMaking.Range("C76").NumberFormat = "@"
$marketMaking.Range("C76").Value = "21:07:32"
$marketMaking.Range("C76").Style = "Normal"
$marketMaking.Range("D76").Value = "MarketMaking"
$marketMaking.Range("E76").Value = "UP"
$marketMaking.Range("F76").Value = 0.86
$marketMaking.Range("H76").Value = "OPEN"
$marketMaking.Range("I76").Value = 0
$marketMaking.Range("J76").Value = 0
$marketMaking.Range("K76").Value = 100.4914872031006
$marketMaking.Range("L76").Value = 0
$marketMaking.Range("M76").Value = 0
$marketMaking.Range("N76").Value = 0.6
$marketMaking.Range("O76").Value = "Normal spread capture: 19600 bps"
$marketMaking.Range("Q76").Value = 0
